$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing blank separator row (row 136).
# This shifts: blank separator 136->138, "sum [min]" 137->139,
# "sum [h]" 138->140, "sum [working weeks]" 139->141.
$ws.Rows("136:137").Insert()

# Copy formatting (number formats / styles) from the last data row (135)
# down into the two freshly inserted rows so D/E/F/G pick up the same
# styles used by the rest of the time-entry rows.
$ws.Range("A135:G135").Copy()
$ws.Range("A136:G137").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New data row 136: 2014-07-12, 23:00 -> 24:00 (59.99999.. minutes / ~1h)
$ws.Cells.Item(136, 1).Value = 2014
$ws.Cells.Item(136, 2).Value = 7
$ws.Cells.Item(136, 3).Value = 12
$ws.Cells.Item(136, 4).Value = 0.95833333333333337
$ws.Cells.Item(136, 5).Value = 1
$ws.Range("F136").Formula = "=(E136-D136)*24*60"
$ws.Range("G136").Formula = "=F136/60"

# New data row 137: 2014-07-13, 00:00 -> 00:00 (zero-length entry)
$ws.Cells.Item(137, 1).Value = 2014
$ws.Cells.Item(137, 2).Value = 7
$ws.Cells.Item(137, 3).Value = 13
$ws.Cells.Item(137, 4).Value = 0
$ws.Cells.Item(137, 5).Value = 0
$ws.Range("F137").Formula = "=(E137-D137)*24*60"
$ws.Range("G137").Formula = "=F137/60"

# Update the SUM range (row 139 now, was 137) to include the new rows.
$ws.Range("F139").Formula = "=SUM(F2:F137)"

# Selection follows the blank separator row, which is now E138.
$ws.Range("E138").Select()

$wb.Application.Calculate()
